$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same cell formatting (fill + wrap text) used elsewhere in column B
# (style index 3: themed fill + wrap text) to B24:B26 by copying format from
# a cell that already carries that style (e.g. B11).
$ws.Range("B11").Copy()
$ws.Range("B24:B26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection shown in the sheet view.
$ws.Range("A25").Select()
